# Adds match-day "N3 J12 ES Fos" (2026-01-09) player rows 1104-1117 to the tracking sheet,
# mirroring the other "Global" match rows already present (rows 2-1103).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstNewRow = 1104
$lastNewRow  = 1117

# Seed the new block by duplicating the last existing data row so the new cells
# inherit the same number formats / styles (date style on B, centred style on D, ...).
$ws.Range("A1103:V1103").Copy($ws.Range("A{0}:V{1}" -f $firstNewRow, $lastNewRow))

# --- Column G (time on pitch) -- filled row by row, first, matching entry order ---
$ws.Range("G1104").Value = "01:24:21"
$ws.Range("G1105").Value = "01:39:26"
$ws.Range("G1106").Value = "00:14:04"
$ws.Range("G1107").Value = "01:39:26"
$ws.Range("G1108").Value = "01:37:48"
$ws.Range("G1109").Value = "01:24:13"
$ws.Range("G1110").Value = "00:05:05"
$ws.Range("G1111").Value = "01:39:19"
$ws.Range("G1112").Value = "01:39:26"
$ws.Range("G1113").Value = "00:37:33"
$ws.Range("G1114").Value = "01:00:53"
$ws.Range("G1115").Value = "01:32:36"
$ws.Range("G1116").Value = "00:12:56"
$ws.Range("G1117").Value = "00:44:24"

# --- Columns E (player name) and F (position) ---
$ws.Range("E1104").Value = "Sofiane Belle"
$ws.Range("F1104").Value = "left forward"
$ws.Range("E1105").Value = "Yoan Zouma"
$ws.Range("F1105").Value = "center back"
$ws.Range("E1106").Value = "Karahali Souaré"
$ws.Range("F1106").Value = "right forward"
$ws.Range("E1107").Value = "Malik Boussaid"
$ws.Range("F1107").Value = "right back"
$ws.Range("E1108").Value = "Maé Clavel"
$ws.Range("F1108").Value = "left back"
$ws.Range("E1109").Value = "Emmanuel Valey"
$ws.Range("F1109").Value = "left forward"
$ws.Range("E1110").Value = "Theo Owono"
$ws.Range("F1110").Value = "center midfield"
$ws.Range("E1111").Value = "Naim Dhib"
$ws.Range("F1111").Value = "center midfield"
$ws.Range("E1112").Value = "Romain Thunet"
$ws.Range("F1112").Value = "center back"
$ws.Range("E1113").Value = "Amir Etien"
$ws.Range("F1113").Value = "right forward"
$ws.Range("E1114").Value = "Ilan Ihaddadene"
$ws.Range("F1114").Value = "center midfield"
$ws.Range("E1115").Value = "Yoann Martelat"
$ws.Range("F1115").Value = "center midfield"
$ws.Range("E1116").Value = "Karim Belmahi"
$ws.Range("F1116").Value = "left forward"
$ws.Range("E1117").Value = "Naim Ighbane"
$ws.Range("F1117").Value = "center back"

# --- Columns B (match date), C (category), D (match marker) ---
$ws.Range("B1104").Value = 46031
$ws.Range("C1104").Value = "Global"
$ws.Range("D1104").Value = "M"
$ws.Range("B1105").Value = 46031
$ws.Range("C1105").Value = "Global"
$ws.Range("D1105").Value = "M"
$ws.Range("B1106").Value = 46031
$ws.Range("C1106").Value = "Global"
$ws.Range("D1106").Value = "M"
$ws.Range("B1107").Value = 46031
$ws.Range("C1107").Value = "Global"
$ws.Range("D1107").Value = "M"
$ws.Range("B1108").Value = 46031
$ws.Range("C1108").Value = "Global"
$ws.Range("D1108").Value = "M"
$ws.Range("B1109").Value = 46031
$ws.Range("C1109").Value = "Global"
$ws.Range("D1109").Value = "M"
$ws.Range("B1110").Value = 46031
$ws.Range("C1110").Value = "Global"
$ws.Range("D1110").Value = "M"
$ws.Range("B1111").Value = 46031
$ws.Range("C1111").Value = "Global"
$ws.Range("D1111").Value = "M"
$ws.Range("B1112").Value = 46031
$ws.Range("C1112").Value = "Global"
$ws.Range("D1112").Value = "M"
$ws.Range("B1113").Value = 46031
$ws.Range("C1113").Value = "Global"
$ws.Range("D1113").Value = "M"
$ws.Range("B1114").Value = 46031
$ws.Range("C1114").Value = "Global"
$ws.Range("D1114").Value = "M"
$ws.Range("B1115").Value = 46031
$ws.Range("C1115").Value = "Global"
$ws.Range("D1115").Value = "M"
$ws.Range("B1116").Value = 46031
$ws.Range("C1116").Value = "Global"
$ws.Range("D1116").Value = "M"
$ws.Range("B1117").Value = 46031
$ws.Range("C1117").Value = "Global"
$ws.Range("D1117").Value = "M"

# --- Numeric performance stat columns H through V ---
$ws.Range("H1104").Value = 8.26
$ws.Range("I1104").Value = 1.4
$ws.Range("J1104").Value = 6.84
$ws.Range("K1104").Value = 0.91
$ws.Range("L1104").Value = 0.41
$ws.Range("M1104").Value = 0.1
$ws.Range("N1104").Value = 0
$ws.Range("O1104").Value = 8
$ws.Range("P1104").Value = 5.86
$ws.Range("Q1104").Value = 28.64
$ws.Range("R1104").Value = 4.67
$ws.Range("S1104").Value = 32
$ws.Range("T1104").Value = 6
$ws.Range("U1104").Value = 29
$ws.Range("V1104").Value = 10
$ws.Range("H1105").Value = 9.32
$ws.Range("I1105").Value = 1.23
$ws.Range("J1105").Value = 8.08
$ws.Range("K1105").Value = 0.81
$ws.Range("L1105").Value = 0.33
$ws.Range("M1105").Value = 0.1
$ws.Range("N1105").Value = 0.01
$ws.Range("O1105").Value = 8
$ws.Range("P1105").Value = 5.51
$ws.Range("Q1105").Value = 30.41
$ws.Range("R1105").Value = 4.47
$ws.Range("S1105").Value = 26
$ws.Range("T1105").Value = 7
$ws.Range("U1105").Value = 19
$ws.Range("V1105").Value = 12
$ws.Range("H1106").Value = 1.84
$ws.Range("I1106").Value = 0.55
$ws.Range("J1106").Value = 1.28
$ws.Range("K1106").Value = 0.35
$ws.Range("L1106").Value = 0.18
$ws.Range("M1106").Value = 0.03
$ws.Range("N1106").Value = 0
$ws.Range("O1106").Value = 2
$ws.Range("P1106").Value = 7.75
$ws.Range("Q1106").Value = 28.29
$ws.Range("R1106").Value = 3.99
$ws.Range("S1106").Value = 12
$ws.Range("T1106").Value = 0
$ws.Range("U1106").Value = 8
$ws.Range("V1106").Value = 2
$ws.Range("H1107").Value = 11.3
$ws.Range("I1107").Value = 2.16
$ws.Range("J1107").Value = 9.12
$ws.Range("K1107").Value = 1.39
$ws.Range("L1107").Value = 0.63
$ws.Range("M1107").Value = 0.16
$ws.Range("N1107").Value = 0
$ws.Range("O1107").Value = 11
$ws.Range("P1107").Value = 6.63
$ws.Range("Q1107").Value = 30.28
$ws.Range("R1107").Value = 4.25
$ws.Range("S1107").Value = 55
$ws.Range("T1107").Value = 7
$ws.Range("U1107").Value = 51
$ws.Range("V1107").Value = 15
$ws.Range("H1108").Value = 11.24
$ws.Range("I1108").Value = 1.76
$ws.Range("J1108").Value = 9.45
$ws.Range("K1108").Value = 1.31
$ws.Range("L1108").Value = 0.39
$ws.Range("M1108").Value = 0.09
$ws.Range("N1108").Value = 0
$ws.Range("O1108").Value = 7
$ws.Range("P1108").Value = 6.82
$ws.Range("Q1108").Value = 28.1
$ws.Range("R1108").Value = 4.47
$ws.Range("S1108").Value = 36
$ws.Range("T1108").Value = 3
$ws.Range("U1108").Value = 32
$ws.Range("V1108").Value = 14
$ws.Range("H1109").Value = 9.56
$ws.Range("I1109").Value = 2.09
$ws.Range("J1109").Value = 7.44
$ws.Range("K1109").Value = 1.31
$ws.Range("L1109").Value = 0.65
$ws.Range("M1109").Value = 0.16
$ws.Range("N1109").Value = 0
$ws.Range("O1109").Value = 16
$ws.Range("P1109").Value = 6.72
$ws.Range("Q1109").Value = 29.65
$ws.Range("R1109").Value = 4.95
$ws.Range("S1109").Value = 50
$ws.Range("T1109").Value = 12
$ws.Range("U1109").Value = 41
$ws.Range("V1109").Value = 18
$ws.Range("H1110").Value = 0.62
$ws.Range("I1110").Value = 0.16
$ws.Range("J1110").Value = 0.45
$ws.Range("K1110").Value = 0.1
$ws.Range("L1110").Value = 0.07
$ws.Range("M1110").Value = 0
$ws.Range("N1110").Value = 0
$ws.Range("O1110").Value = 1
$ws.Range("P1110").Value = 6.85
$ws.Range("Q1110").Value = 25.17
$ws.Range("R1110").Value = 4.47
$ws.Range("S1110").Value = 2
$ws.Range("T1110").Value = 2
$ws.Range("U1110").Value = 2
$ws.Range("V1110").Value = 1
$ws.Range("H1111").Value = 10.43
$ws.Range("I1111").Value = 1.59
$ws.Range("J1111").Value = 8.82
$ws.Range("K1111").Value = 1.22
$ws.Range("L1111").Value = 0.34
$ws.Range("M1111").Value = 0.04
$ws.Range("N1111").Value = 0.01
$ws.Range("O1111").Value = 3
$ws.Range("P1111").Value = 6.23
$ws.Range("Q1111").Value = 31.17
$ws.Range("R1111").Value = 4.42
$ws.Range("S1111").Value = 60
$ws.Range("T1111").Value = 4
$ws.Range("U1111").Value = 37
$ws.Range("V1111").Value = 9
$ws.Range("H1112").Value = 10.53
$ws.Range("I1112").Value = 1.27
$ws.Range("J1112").Value = 9.24
$ws.Range("K1112").Value = 0.87
$ws.Range("L1112").Value = 0.27
$ws.Range("M1112").Value = 0.14
$ws.Range("N1112").Value = 0.01
$ws.Range("O1112").Value = 7
$ws.Range("P1112").Value = 6.29
$ws.Range("Q1112").Value = 30.4
$ws.Range("R1112").Value = 4.46
$ws.Range("S1112").Value = 30
$ws.Range("T1112").Value = 5
$ws.Range("U1112").Value = 28
$ws.Range("V1112").Value = 12
$ws.Range("H1113").Value = 3.41
$ws.Range("I1113").Value = 0.62
$ws.Range("J1113").Value = 2.78
$ws.Range("K1113").Value = 0.31
$ws.Range("L1113").Value = 0.24
$ws.Range("M1113").Value = 0.08
$ws.Range("N1113").Value = 0
$ws.Range("O1113").Value = 6
$ws.Range("P1113").Value = 5.29
$ws.Range("Q1113").Value = 30.17
$ws.Range("R1113").Value = 4.94
$ws.Range("S1113").Value = 19
$ws.Range("T1113").Value = 4
$ws.Range("U1113").Value = 14
$ws.Range("V1113").Value = 6
$ws.Range("H1114").Value = 7.1
$ws.Range("I1114").Value = 1.37
$ws.Range("J1114").Value = 5.71
$ws.Range("K1114").Value = 1.06
$ws.Range("L1114").Value = 0.31
$ws.Range("M1114").Value = 0.02
$ws.Range("N1114").Value = 0
$ws.Range("O1114").Value = 3
$ws.Range("P1114").Value = 6.83
$ws.Range("Q1114").Value = 26.9
$ws.Range("R1114").Value = 4.64
$ws.Range("S1114").Value = 26
$ws.Range("T1114").Value = 10
$ws.Range("U1114").Value = 29
$ws.Range("V1114").Value = 6
$ws.Range("H1115").Value = 11.1
$ws.Range("I1115").Value = 2.2
$ws.Range("J1115").Value = 8.88
$ws.Range("K1115").Value = 1.79
$ws.Range("L1115").Value = 0.39
$ws.Range("M1115").Value = 0.04
$ws.Range("N1115").Value = 0
$ws.Range("O1115").Value = 4
$ws.Range("P1115").Value = 7.17
$ws.Range("Q1115").Value = 28.84
$ws.Range("R1115").Value = 4.35
$ws.Range("S1115").Value = 30
$ws.Range("T1115").Value = 6
$ws.Range("U1115").Value = 24
$ws.Range("V1115").Value = 13
$ws.Range("H1116").Value = 1.75
$ws.Range("I1116").Value = 0.5
$ws.Range("J1116").Value = 1.24
$ws.Range("K1116").Value = 0.36
$ws.Range("L1116").Value = 0.14
$ws.Range("M1116").Value = 0.01
$ws.Range("N1116").Value = 0
$ws.Range("O1116").Value = 2
$ws.Range("P1116").Value = 8.05
$ws.Range("Q1116").Value = 25.52
$ws.Range("R1116").Value = 4.03
$ws.Range("S1116").Value = 10
$ws.Range("T1116").Value = 1
$ws.Range("U1116").Value = 7
$ws.Range("V1116").Value = 6
$ws.Range("H1117").Value = 4.65
$ws.Range("I1117").Value = 0.69
$ws.Range("J1117").Value = 3.95
$ws.Range("K1117").Value = 0.42
$ws.Range("L1117").Value = 0.21
$ws.Range("M1117").Value = 0.06
$ws.Range("N1117").Value = 0.01
$ws.Range("O1117").Value = 5
$ws.Range("P1117").Value = 6.16
$ws.Range("Q1117").Value = 30.6
$ws.Range("R1117").Value = 4.43
$ws.Range("S1117").Value = 17
$ws.Range("T1117").Value = 3
$ws.Range("U1117").Value = 17
$ws.Range("V1117").Value = 7

# --- Column A (match label) -- set last, across the whole new block in one shot, ---
# --- matching how the source workbook only records a single new shared string for it ---
$ws.Range("A{0}:A{1}" -f $firstNewRow, $lastNewRow).Value = "N3 J12 ES Fos"

# Leave the selection where the author left it after entering the new match.
[void]$ws.Range("F1120").Select()